$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.194.69'
$ws.Range("E2").Value = '  +9.17%  '

$ws.Range("D3").Value = '3.239.32'
$ws.Range("E3").Value = '  +4.38%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '401.07'
$ws.Range("E5").Value = '  +4.43%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '111.31'
$ws.Range("E6").Value = '  +7.94%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.557'
$ws.Range("E7").Value = '  +3.18%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.628'
$ws.Range("E9").Value = '  +7.67%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.68'
$ws.Range("E10").Value = '  +6.83%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0902'
$ws.Range("E11").Value = '  +5.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.141'
$ws.Range("E12").Value = '  +2.16%  '

$ws.Range("D13").Value = '3.743.68'
$ws.Range("E13").Value = '  +4.16%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.21'
$ws.Range("E14").Value = '  +3.27%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.11'
$ws.Range("E15").Value = '  +3.81%  '

$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.07'
$ws.Range("E16").Value = '  +7.78%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.232.73'
$ws.Range("E17").Value = '  +4.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '10.57'
$ws.Range("E18").Value = '  -7.19%  '

$ws.Range("D19").Value = '56.239.53'
$ws.Range("E19").Value = '  +9.23%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.38'
$ws.Range("E20").Value = '  +2.70%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.20'
$ws.Range("E21").Value = '  +7.14%  '

$ws.Range("E22").Value = '  +5.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '289.25'
$ws.Range("E23").Value = '  +8.87%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.71'
$ws.Range("E24").Value = '  +6.85%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.25'
$ws.Range("E25").Value = '  +4.83%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.22'
$ws.Range("E26").Value = '  +1.63%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '28.23'
$ws.Range("E27").Value = '  +4.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.58'
$ws.Range("E28").Value = '  +4.34%  '

$ws.Range("E29").Value = '  +3.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.112'
$ws.Range("E31").Value = '  +4.69%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.35'
$ws.Range("E32").Value = '  +10.05%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0500'
$ws.Range("E33").Value = '  +6.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '37.21'
$ws.Range("E34").Value = '  +5.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.11'
$ws.Range("E35").Value = '  +1.85%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.21'
$ws.Range("E36").Value = '  +1.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.59'
$ws.Range("E37").Value = '  +7.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.08%  '

$ws.Range("E39").Value = '  +21.91%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '138.57'
$ws.Range("E40").Value = '  +7.40%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.94'
$ws.Range("E41").Value = '  +2.84%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.07'
$ws.Range("E42").Value = '  +10.82%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.288'
$ws.Range("E43").Value = '  -3.03%  '

$ws.Range("E44").Value = '  +2.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.83'
$ws.Range("E45").Value = '  +1.95%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.63'
$ws.Range("E46").Value = '  +2.12%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.12'
$ws.Range("E47").Value = '  +42.18%  '

$ws.Range("E48").Value = '  +0.22%  '

$ws.Range("B49").Value = 'Maker'
$ws.Range("C49").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D49").Value = '2.135.37'
$ws.Range("E49").Value = '  +3.87%  '

$ws.Range("B50").Value = 'WEMIXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.09'
$ws.Range("E50").Value = '  +0.85%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0346'
$ws.Range("E51").Value = '  +6.52%  '
